# Apply updated Betfair Back/Lay odds for 2025-11-14 to row 3..8
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 (HIK Hellerup vs Vendsyssel FF) ---
$ws.Range("G3").Value = 3.65
$ws.Range("I3").Value = 2.6
$ws.Range("P3").Value = 1.82
$ws.Range("W3").Value = 1.37
$ws.Range("Z3").Value = 16.5
$ws.Range("AA3").Value = 40
$ws.Range("AD3").Value = 14
$ws.Range("AG3").Value = 15
$ws.Range("AL3").Value = 1000

# --- Row 4 (Valladolid vs Las Palmas) ---
$ws.Range("F4").Value = 2.36
$ws.Range("G4").Value = 2.5
$ws.Range("J4").Value = 3.05
$ws.Range("L4").Value = 1.56
$ws.Range("N4").Value = 2.68
$ws.Range("O4").Value = 1.52
$ws.Range("P4").Value = 1.56
$ws.Range("Q4").Value = 2.28
$ws.Range("W4").Value = 1.66
$ws.Range("X4").Value = 8.800000000000001
$ws.Range("AB4").Value = 7.8
$ws.Range("AG4").Value = 12.5
$ws.Range("AH4").Value = 23
$ws.Range("AK4").Value = 34
$ws.Range("AM4").Value = 190
$ws.Range("AN4").Value = 36

# --- Row 5 (Flint Town United vs The New Saints) ---
$ws.Range("H5").Value = 1.23
$ws.Range("I5").Value = 1.3
$ws.Range("L5").Value = 1.01
$ws.Range("N5").Value = 8
$ws.Range("P5").Value = 3.35
$ws.Range("R5").Value = 1.97
$ws.Range("T5").Value = 1.77
$ws.Range("V5").Value = 4.3
$ws.Range("W5").Value = 1.07

# --- Row 6 (Cardiff Metropolitan vs Briton Ferry Llansawel) ---
$ws.Range("I6").Value = 5.6
$ws.Range("P6").Value = 2.14
$ws.Range("U6").Value = 2.12
$ws.Range("X6").Value = 20
$ws.Range("Y6").Value = 21

# --- Row 7 (Lanus vs Atl Tucuman) ---
$ws.Range("F7").Value = 1.92
$ws.Range("K7").Value = 3.85

# --- Row 8 (Paysandu vs Amazonas FC) ---
$ws.Range("F8").Value = 3.2
$ws.Range("G8").Value = 3.6
$ws.Range("H8").Value = 2.38
$ws.Range("I8").Value = 2.56
$ws.Range("J8").Value = 3.25
$ws.Range("K8").Value = 3.5
$ws.Range("S8").Value = 4.1
$ws.Range("V8").Value = 1.64
$ws.Range("W8").Value = 1.38
$ws.Range("AA8").Value = 42
$ws.Range("AI8").Value = 60
$ws.Range("AJ8").Value = 70
$ws.Range("AK8").Value = 60
$ws.Range("AM8").Value = 150
$ws.Range("AN8").Value = 1000
